# Update the "department" column (C) values from "FACULTY OF TECH SCIENCES"
# to "Automotive" for the course rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Automotive"
$ws.Range("C3").Value = "Automotive"
$ws.Range("C4").Value = "Automotive"
